$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# --- New journal entry on row 3 ---
# Pull formatting (number formats / styles) from row 2's cells first so the
# new cells inherit the same look (date format on A, time format on F:H).
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("F2:H2").Copy()
$ws.Range("F3:H3").PasteSpecial(-4122)

$ws.Range("A3").Value = 44958
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Documentation"
$ws.Range("D3").Value = "Recherche documentation"
$ws.Range("E3").Value = "Normes ISO:27001 afin de créer un politique de sécurité de l'information pour mon entreprise"
$ws.Range("F3").Value = 0.44444444444444442
$ws.Range("G3").Value = 0.51041666666666663
$ws.Range("H3").Formula = "=G3-F3"

# --- Pre-format / pre-compute the rest of the log rows (4-50), mirroring
#     the blank template rows that already carried the date/time formats ---
$ws.Range("A2").Copy()
$ws.Range("A4:A36").PasteSpecial(-4122)
$ws.Range("F2:H2").Copy()
$ws.Range("F4:H50").PasteSpecial(-4122)
$ws.Range("H4:H50").FormulaR1C1 = "=RC[-1]-RC[-2]"

# --- Total row ---
$ws.Range("H2").Copy()
$ws.Range("H51").PasteSpecial(-4122)
$ws.Range("H51").Formula = "=SUM(H2:H50)"

$excel.CutCopyMode = 0

# --- Restore the selection to where the user last clicked ---
[void]$ws.Range("H52").Select()
